$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row
# (rows 2 through 288). Bump every one of these dates forward by one day:
# serial 45181 (2023-09-12) -> serial 45182 (2023-09-13).
$startRow = 2
$endRow = 288

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -eq 45181) {
        $cell.Value = 45182
    }
}
